# Update bat 79 brain regions of each TT (tetrode) for all days.
# Rows 31-55 on the "Experiments" sheet correspond to bat 79's recording
# days. Column K ("TT_loc") gets a per-tetrode brain-region label, and
# columns N/O (neural_data_exist / position_data_exist) get flagged to 1
# for every row that didn't already carry that flag.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row -> TT_loc (column K) text value.
$ttLoc = @{
    31 = "{'WM','WM','WM','N/A'}"
    32 = "{'WM','WM','WM','N/A'}"
    33 = "{'WM','WM','CA1','N/A'}"
    34 = "{'WM','WM','CA1','N/A'}"
    35 = "{'WM','WM','CA1','N/A'}"
    36 = "{'WM','WM','CA1','N/A'}"
    37 = "{'WM','CA1','CA1','N/A'}"
    38 = "{'WM','CA1','CA1','N/A'}"
    39 = "{'WM','CA1','CA1','N/A'}"
    40 = "{'WM','CA1','CA1','N/A'}"
    41 = "{'WM','CA1','CA1','N/A'}"
    42 = "{'WM','CA1','CA1','N/A'}"
    43 = "{'WM','CA1','CA1','N/A'}"
    44 = "{'WM','CA1','CA1','N/A'}"
    45 = "{'WM','CA1','CA1','N/A'}"
    46 = "{'WM','CA1','CA1','N/A'}"
    47 = "{'WM','CA1','CA1','N/A'}"
    48 = "{'CA1','CA1','CA1','N/A'}"
    49 = "{'CA1','CA1','CA1','N/A'}"
    50 = "{'CA1','CA1','CA1','N/A'}"
    51 = "{'CA1','CA1','CA1','N/A'}"
    52 = "{'CA1','CA1','CA1','N/A'}"
    53 = "{'CA1','CA1','CA1','N/A'}"
    54 = "{'CA1','CA1','CA1','N/A'}"
    55 = "{'CA1','CA1','CA1','N/A'}"
}

# Rows that already had neural_data_exist (N) / position_data_exist (O)
# populated before this edit - leave those alone.
$alreadyFlagged = @(37, 42, 48, 50, 51)

foreach ($row in 31..55) {
    $ws.Cells.Item($row, 11).Value = $ttLoc[$row]   # column K

    if ($alreadyFlagged -notcontains $row) {
        $ws.Cells.Item($row, 14).Value = 1          # column N
        $ws.Cells.Item($row, 15).Value = 1          # column O
    }
}

# Update the sheet's active selection (bottom-right frozen pane) to match
# the author's final cursor position.
$ws.Range("K49").Select() | Out-Null
